$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E5").Style = "Normal"
$ws.Range("E5").Value2 = "yes"
